$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit "fixed for some error configure file" corrects the bad value
# in column M ("CanClone") for the SelectScene/City row (row 6): it was
# mistakenly 0 and should be 1.
$ws.Range("M6").Value = 1

# Reflect the cell that was active/selected in the sheet when the author
# saved the corrected workbook.
$ws.Range("N12").Select() | Out-Null
